$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 22-29 replicate rows 8-15 (second/third cycle block repeated as a new "cuarto ciclo" block)
$ws.Range("A8:P15").Copy($ws.Range("A22:P29"))

# Row 30 is a section separator row with only the "CUARTO CICLO" label in column B
$ws.Range("A21:P21").Copy($ws.Range("A30:P30"))
$ws.Range("A30:P30").ClearContents()
$ws.Cells.Item(30, 2).Value = "CUARTO CICLO"

# Rows 31-32 replicate rows 8-9 again
$ws.Range("A8:P9").Copy($ws.Range("A31:P32"))

# Restore the active selection to match the new last cell used (P32)
$ws.Range("P32").Select()
